$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder "Materias primas" (raw materials) text for a few products so that
# Vainilla / Limon entries move to the end of the comma-separated list.
$ws.Range("C2").Value = "1.0-Huevos (unidad),3.0-Leche (litros),2.0-Harina  (kg),1.0-Vainilla (ml),"
$ws.Range("C4").Value = "2.0-Huevos (unidad),5.0-Harina  (kg),1.0-Vainilla (ml),"
$ws.Range("C6").Value = "5.0-Huevos (unidad),5.0-Crema (litros),4.0-Harina  (kg),2.0-Limon (unidad),"
$ws.Range("C7").Value = "2.0-Huevos (unidad),0.2-Leche (litros),0.3-Harina  (kg),0.1-Vainilla (ml),"

# "Disponible" column: flip availability flag from 0 to 1 for all products
# (they all share the same underlying string value). Force text storage so
# the value round-trips as a shared string instead of a number.
$ws.Range("E2:E7").NumberFormat = "@"
$ws.Range("E2:E7").Value = "1"
